{"js": "// Apply the author's edits to the GZJ-zlecenie contract:\n//   * The contracting individual changes from a woman (Pani Katarzyna Faller)\n//     to a man (Pan Maciej Karli\u0144ski) \u2014 name, surname, honorific and the\n//     grammatical gender suffixes (\"...\u0105\" -> \"...ym\") all change in lockstep.\n//   * Their address, PESEL number and bank account number are replaced with\n//     new values.\n//\n// Every search string below is a unique, verbatim run of text in the\n// document, so Body.search(..., { matchCase: true }) followed by\n// Range.insertText(..., \"Replace\") on the single hit is a safe, unambiguous\n// way to perform each substitution.\n\nconst body = context.document.body;\n\nasync function replaceUniqueText(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for '\" + findText + \"' but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// Honorific: \"Pani\u0105\" (Mrs., accusative) -> \"Panem\" (Mr., instrumental)\nawait replaceUniqueText(\"Pani\u0105\", \"Panem\");\n\n// First name: Katarzyn\u0105 -> Maciejem\nawait replaceUniqueText(\"Katarzyn\u0105\", \"Maciejem\");\n\n// Surname: Faller -> Karli\u0144skim\nawait replaceUniqueText(\"Faller\", \"Karli\u0144skim\");\n\n// Postal code: 61-131 -> 61-211\nawait replaceUniqueText(\"61-131\", \"61-211\");\n\n// Street: ul. Katowicka -> os. O\u015bwiecenia\nawait replaceUniqueText(\"ul. Katowicka\", \"os. O\u015bwiecenia\");\n\n// House/flat number: 53A/45 -> 98/60\nawait replaceUniqueText(\"53A/45\", \"98/60\");\n\n// PESEL number\nawait replaceUniqueText(\"93092305326\", \"92062414370\");\n\n// Grammatical-gender agreement (feminine \"-\u0105\" -> masculine \"-ym\") for the\n// three participle/adjective forms describing the contractor.\nawait replaceUniqueText(\"reprezentuj\u0105c\u0105\", \"reprezentuj\u0105cym\");\nawait replaceUniqueText(\"zapewniaj\u0105c\u0105\", \"zapewniaj\u0105cym\");\nawait replaceUniqueText(\"zwan\u0105\", \"zwanym\");\n\n// Bank account number\nawait replaceUniqueText(\n  \"84 1140 2004 0000 3202 7942 7526\",\n  \"28 1140 2004 0000 3402 7054 1393\"\n);\n", "ps1": "# Apply the author's edits to the GZJ-zlecenie contract:\n#   * The contracting individual changes from a woman (Pani Katarzyna Faller)\n#     to a man (Pan Maciej Karlinski) - name, surname, honorific and the\n#     grammatical gender suffixes (\"...a\" -> \"...ym\") all change in lockstep.\n#   * Their address, PESEL number and bank account number are replaced with\n#     new values.\n#\n# Every search string below is a unique, verbatim run of text in the\n# document, so a case-sensitive Find/Replace-All (which, given uniqueness,\n# touches exactly one spot) is a safe way to perform each substitution.\n\n$d = $word.ActiveDocument\n\nfunction Replace-UniqueText($findText, $replaceText) {\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $findText,    # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        $true,        # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Could not find expected text: '$findText'\"\n    }\n}\n\n# Honorific: \"Pani\u0105\" (Mrs., accusative) -> \"Panem\" (Mr., instrumental)\nReplace-UniqueText \"Pani\u0105\" \"Panem\"\n\n# First name: Katarzyn\u0105 -> Maciejem\nReplace-UniqueText \"Katarzyn\u0105\" \"Maciejem\"\n\n# Surname: Faller -> Karli\u0144skim\nReplace-UniqueText \"Faller\" \"Karli\u0144skim\"\n\n# Postal code: 61-131 -> 61-211\nReplace-UniqueText \"61-131\" \"61-211\"\n\n# Street: ul. Katowicka -> os. O\u015bwiecenia\nReplace-UniqueText \"ul. Katowicka\" \"os. O\u015bwiecenia\"\n\n# House/flat number: 53A/45 -> 98/60\nReplace-UniqueText \"53A/45\" \"98/60\"\n\n# PESEL number\nReplace-UniqueText \"93092305326\" \"92062414370\"\n\n# Grammatical-gender agreement (feminine \"-\u0105\" -> masculine \"-ym\") for the\n# three participle/adjective forms describing the contractor.\nReplace-UniqueText \"reprezentuj\u0105c\u0105\" \"reprezentuj\u0105cym\"\nReplace-UniqueText \"zapewniaj\u0105c\u0105\" \"zapewniaj\u0105cym\"\nReplace-UniqueText \"zwan\u0105\" \"zwanym\"\n\n# Bank account number\nReplace-UniqueText \"84 1140 2004 0000 3202 7942 7526\" \"28 1140 2004 0000 3402 7054 1393\"\n"}
